$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Total number of Patents" (E) and "Total number of Citations" (F)
# for the first two data rows.
$ws.Range("E2").Value = 187.0
$ws.Range("F2").Value = 61.0

$ws.Range("E3").Value = 112.0
$ws.Range("F3").Value = 24.0
